$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set row labels in the order the strings were first authored (matches shared-string table order)
$ws.Cells.Item(3, 1).Value = "volume_adi"
$ws.Cells.Item(4, 1).Value = "volume_obv"
$ws.Cells.Item(5, 1).Value = "volume_cmf"
$ws.Cells.Item(6, 1).Value = "volume_fi"
$ws.Cells.Item(7, 1).Value = "volume_mfi"
$ws.Cells.Item(8, 1).Value = "volume_em"
$ws.Cells.Item(9, 1).Value = "volume_sma_em"
$ws.Cells.Item(10, 1).Value = "volume_vpt"
$ws.Cells.Item(11, 1).Value = "volume_nvi"
$ws.Cells.Item(12, 1).Value = "volume_vwap"
$ws.Cells.Item(13, 1).Value = "volatility_atr"
$ws.Cells.Item(14, 1).Value = "volatility_bbm"
$ws.Cells.Item(15, 1).Value = "volatility_bbl"
$ws.Cells.Item(16, 1).Value = "volatility_bbw"
$ws.Cells.Item(17, 1).Value = "volatility_bbp"
$ws.Cells.Item(18, 1).Value = "volatility_bbhi"
$ws.Cells.Item(19, 1).Value = "volatility_bbli"
$ws.Cells.Item(20, 1).Value = "volatility_kch"
$ws.Cells.Item(21, 1).Value = "volatility_kcl"
$ws.Cells.Item(22, 1).Value = "volatility_kcw"
$ws.Cells.Item(23, 1).Value = "volatility_kcp"
$ws.Cells.Item(24, 1).Value = "volatility_kchi"
$ws.Cells.Item(25, 1).Value = "volatility_kcli"
$ws.Cells.Item(2, 1).Value = "baseline_norm"
$ws.Cells.Item(1, 1).Value = "baseline"

# Set numeric values (row order, ascending)
$ws.Cells.Item(1, 2).Value = 1.3369819999999999
$ws.Cells.Item(2, 2).Value = 1.5157309999999999
$ws.Cells.Item(3, 2).Value = 1.2258605593398599
$ws.Cells.Item(4, 2).Value = 1.8479266536594201
$ws.Cells.Item(5, 2).Value = 1.2823241926381801
$ws.Cells.Item(6, 2).Value = 1.17036116564987
$ws.Cells.Item(7, 2).Value = 1.2288814881314101
$ws.Cells.Item(8, 2).Value = 1.28509253366123
$ws.Cells.Item(9, 2).Value = 1.1045496052777799
$ws.Cells.Item(10, 2).Value = 1.11740781881705
$ws.Cells.Item(11, 2).Value = 1.21779312217193
$ws.Cells.Item(12, 2).Value = 1.20205461919401
$ws.Cells.Item(13, 2).Value = 1.0824852489251899
$ws.Cells.Item(14, 2).Value = 1.62919189390126
$ws.Cells.Item(15, 2).Value = 1.2788475072946299
$ws.Cells.Item(16, 2).Value = 1.2766168844767101
$ws.Cells.Item(17, 2).Value = 1.4384576680093899
$ws.Cells.Item(18, 2).Value = 1.06792475035613
$ws.Cells.Item(19, 2).Value = 1.14880808000812
$ws.Cells.Item(20, 2).Value = 1.35366429649011
$ws.Cells.Item(21, 2).Value = 1.3992043617385199
$ws.Cells.Item(22, 2).Value = 1.0349786215915999
$ws.Cells.Item(23, 2).Value = 1.0531606524371999
$ws.Cells.Item(24, 2).Value = 1.2765007434275899
$ws.Cells.Item(25, 2).Value = 1.1519725833030601

$ws.Range("A14").Select()
